$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.866.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.56%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.761.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +3.54%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'620.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +4.29%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'177.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.40%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.758.54"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +3.46%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.10%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.26%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +3.62%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -5.43%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -1.85%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'40.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.71%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +1.65%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.384.49"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +3.43%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.761.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.28%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'69.922.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.42%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  +0.21%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.70%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'508.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.82%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'16.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.40%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'9.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +5.11%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.724"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.34%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.65%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'86.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.52%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'13.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.67%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'11.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.17%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0000136"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +23.47%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.22%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'2.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.36%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'2.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +5.11%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -3.65%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'31.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.32%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -1.40%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -0.09%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +5.77%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'6.16"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.63%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.336"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.41%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +2.18%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -1.89%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'50.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.66%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'45.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.60%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'423.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.08%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -1.42%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'3.015.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -3.37%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -1.03%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -1.79%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'27.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -3.89%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D50").Value = "'138.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.40%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.13%  "
$ws.Range("E51").Style = "Normal"
